$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-format the existing data block (rows 4-6, columns D-O) from "0.00" to
# a custom "0.0" number format, matching the number-format cleanup that
# accompanied the new 2022 column.
$ws.Range("D4:O6").NumberFormat = "0.0"

# Add the 2022 column (P) mirroring the existing 2021 column (O): same
# per-row formatting (font, borders, alignment, number format), new data.
$ws.Range("O2").Copy($ws.Range("P2"))
$ws.Range("O3").Copy($ws.Range("P3"))
$ws.Range("O4").Copy($ws.Range("P4"))
$ws.Range("O5").Copy($ws.Range("P5"))
$ws.Range("O6").Copy($ws.Range("P6"))

$ws.Range("P3").Value = 2022
$ws.Range("P4").Formula = "=P5/P6*1000"
$ws.Range("P5").Value = 1339.6
$ws.Range("P6").Value = 6300.5

$ws.Range("S4").Select()
